$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2-6 down to 3-7.
# CopyOrigin = xlFormatFromRightOrBelow (1) so the new row doesn't inherit the
# bold/bordered header formatting from row 1.
$ws.Rows.Item(2).Insert(-4121, 1)

# Clear any inherited formatting on the new row so plain data cells stay unstyled,
# matching the other data rows (A-C, E-R have no explicit style).
$ws.Range("A2:R2").ClearFormats()

# Copy the style (including number format) of the date cell from row 3 (old row 2)
# into the newly inserted row 2's date cell, so the date formatting matches the other rows.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new row 2 with the new record's values.
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(2, 3).Value = "Los Lagos"
$ws.Cells.Item(2, 4).Value = 44473
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 100112012
$ws.Cells.Item(2, 7).Value = "Espinaca"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(2, 11).Value = 11000
$ws.Cells.Item(2, 12).Value = 11000
$ws.Cells.Item(2, 13).Value = 11000
$ws.Cells.Item(2, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(2, 15).Value = "Región Metropolitana"
$ws.Cells.Item(2, 16).Value = 1100
$ws.Cells.Item(2, 17).Value = 10
$ws.Cells.Item(2, 18).Value = "Hortaliza"
